$d = $word.ActiveDocument

# Mapping of old "NN÷N=" text -> new "NN÷N=" text, applied as
# whole-document Find/Replace (each old token is unique in the doc).
$pairs = @(
    @("76÷6=", "26÷6="),
    @("97÷7=", "80÷7="),
    @("97÷4=", "59÷7="),
    @("67÷4=", "72÷9="),
    @("30÷5=", "50÷2="),
    @("10÷2=", "19÷3="),
    @("21÷2=", "55÷2="),
    @("57÷8=", "14÷2="),
    @("20÷8=", "98÷5="),
    @("98÷3=", "71÷2="),
    @("75÷7=", "66÷9="),
    @("51÷4=", "24÷3="),
    @("27÷3=", "43÷3="),
    @("58÷8=", "83÷4="),
    @("78÷6=", "80÷2="),
    @("77÷6=", "48÷6="),
    @("66÷2=", "40÷9="),
    @("71÷7=", "97÷8="),
    @("61÷2=", "15÷3="),
    @("65÷7=", "48÷3="),
    @("11÷3=", "80÷9="),
    @("30÷2=", "55÷2="),
    @("60÷2=", "50÷6="),
    @("61÷4=", "80÷4="),
    @("36÷2=", "80÷8="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done replacing $($pairs.Count) division expressions."
